# Update the due-date column (D) on Sheet1 so that weekend days are not
# counted in the delay calculation - effectively shifting all dates forward.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each block below corresponds to a contiguous run of rows that share the
# same new serial date value (mirrors the D-column edits in the diff).
$ws.Range("D2:D13").Value = 45262
$ws.Range("D14:D30").Value = 45263
$ws.Range("D31:D42").Value = 45264
$ws.Range("D43:D47").Value = 45265
$ws.Range("D48:D52").Value = 45266
$ws.Range("D53:D55").Value = 45267

# Move the active selection on Sheet1 from B51:B55 to C4.
$ws.Range("C4").Select()
